# Weekly driver report update for 2025-04-20
# Re-rank the "Good Drivers" table (rows 13-18) by Total Samples (column B,
# descending) and refresh the per-driver sample counts / vintage dates that
# came in with this week's data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering + refreshed figures for the Good Drivers block (A13:E18)
# Driver Vintage (column E) values are plain date-stamp text, exactly as
# the source report already stores them for rows 15-18, so we force them
# in as text (leading apostrophe) rather than letting Excel reinterpret
# them as real dates.

$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B13").Value = 445055
$ws.Range("D13").Value = 99.90000000000001
$ws.Range("E13").Value = "'2024-11-10"

$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B14").Value = 77849
$ws.Range("D14").Value = 99.90000000000001
$ws.Range("E14").Value = "'2021-08-18"

$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B15").Value = 34244
$ws.Range("D15").Value = 100
$ws.Range("E15").Value = "'2021-04-27"

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B16").Value = 59673
$ws.Range("D16").Value = 100
$ws.Range("E16").Value = "'2020-08-05"

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B17").Value = 113652
$ws.Range("D17").Value = 100
$ws.Range("E17").Value = "'2020-01-06"

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B18").Value = 56018
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = "'2019-12-14"

Write-Output "Driver summary refreshed for 2025-04-20"
